$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns for the data rows so that
# numeric-looking strings (e.g. "1.002") remain stored as text, matching
# the original inline-string cell contents.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "30.567.17"
$ws.Range("E2").Value = "  -0.04%  "

# Row 3
$ws.Range("D3").Value = "1.915.21"
$ws.Range("E3").Value = "  -0.18%  "

# Row 4
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "244.42"
$ws.Range("E5").Value = "  -0.71%  "

# Row 6
$ws.Range("D6").Value = "1.002"

# Row 7
$ws.Range("D7").Value = "0.4866"
$ws.Range("E7").Value = "  +2.72%  "

# Row 8
$ws.Range("D8").Value = "0.2900"
$ws.Range("E8").Value = "  +0.47%  "

# Row 9
$ws.Range("D9").Value = "0.06733"
$ws.Range("E9").Value = "  -1.40%  "

# Row 10
$ws.Range("D10").Value = "111.17"
$ws.Range("E10").Value = "  +5.73%  "

# Row 11
$ws.Range("D11").Value = "19.41"
$ws.Range("E11").Value = "  +5.90%  "

# Row 12
$ws.Range("D12").Value = "1.917.05"
$ws.Range("E12").Value = "  -0.06%  "

# Row 13
$ws.Range("D13").Value = "0.07567"
$ws.Range("E13").Value = "  -1.57%  "

# Row 14
$ws.Range("D14").Value = "5.356"
$ws.Range("E14").Value = "  +1.75%  "

# Row 15
$ws.Range("D15").Value = "0.6722"

# Row 16
$ws.Range("D16").Value = "294.32"
$ws.Range("E16").Value = "  +0.39%  "

# Row 17
$ws.Range("D17").Value = "30.585.61"
$ws.Range("E17").Value = "  +0.01%  "

# Row 18
$ws.Range("D18").Value = "13.04"
$ws.Range("E18").Value = "  +0.94%  "

# Row 19
$ws.Range("E19").Value = "  +0.13%  "

# Row 20
$ws.Range("D20").Value = "0.000007559"
$ws.Range("E20").Value = "  -0.46%  "

# Row 21
$ws.Range("D21").Value = "2.181.04"
$ws.Range("E21").Value = "  +0.73%  "

# Row 22
$ws.Range("D22").Value = "5.504"
$ws.Range("E22").Value = "  +0.31%  "

# Row 23
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.04%  "

# Row 24
$ws.Range("D24").Value = "6.422"
$ws.Range("E24").Value = "  +1.49%  "

# Row 25
$ws.Range("D25").Value = "9.469"
$ws.Range("E25").Value = "  +0.96%  "

# Row 26
$ws.Range("D26").Value = "164.68"
$ws.Range("E26").Value = "  -2.18%  "

# Row 27
$ws.Range("D27").Value = "20.31"
$ws.Range("E27").Value = "  -3.42%  "

# Row 28
$ws.Range("D28").Value = "2.099"
$ws.Range("E28").Value = "  -0.92%  "

# Row 29
$ws.Range("D29").Value = "0.1070"
$ws.Range("E29").Value = "  +0.60%  "

# Row 30
$ws.Range("E30").Value = "  +2.81%  "

# Row 31
$ws.Range("E31").Value = "  -0.51%  "

# Row 32
$ws.Range("D32").Value = "4.063"
$ws.Range("E32").Value = "  -0.17%  "

# Row 33
$ws.Range("D33").Value = "0.04995"
$ws.Range("E33").Value = "  -0.61%  "

# Row 34
$ws.Range("D34").Value = "0.7385"
$ws.Range("E34").Value = "  +0.47%  "

# Row 35
$ws.Range("D35").Value = "1.137"
$ws.Range("E35").Value = "  -0.52%  "

# Row 36
$ws.Range("D36").Value = "0.9998"
$ws.Range("E36").Value = "  +0.05%  "

# Row 37
$ws.Range("D37").Value = "2.714"
$ws.Range("E37").Value = "  -1.17%  "

# Row 38
$ws.Range("E38").Value = "  -1.36%  "

# Row 39
$ws.Range("D39").Value = "2.682"
$ws.Range("E39").Value = "  -0.26%  "

# Row 40
$ws.Range("E40").Value = "  -1.87%  "

# Row 41
$ws.Range("D41").Value = "109.67"
$ws.Range("E41").Value = "  -1.33%  "

# Row 42
$ws.Range("D42").Value = "0.4447"
$ws.Range("E42").Value = "  +1.35%  "

# Row 43
$ws.Range("D43").Value = "0.8625"
$ws.Range("E43").Value = "  -1.68%  "

# Row 44
$ws.Range("D44").Value = "5.842"
$ws.Range("E44").Value = "  -0.56%  "

# Row 45
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "1.002"
$ws.Range("E45").Value = "  +0.18%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "69.67"
$ws.Range("E46").Value = "  +4.06%  "

# Row 47
$ws.Range("D47").Value = "7.231"
$ws.Range("E47").Value = "  -0.25%  "

# Row 48
$ws.Range("D48").Value = "9.245"
$ws.Range("E48").Value = "  -0.21%  "

# Row 49
$ws.Range("D49").Value = "47.94"
$ws.Range("E49").Value = "  +0.20%  "

# Row 50
$ws.Range("D50").Value = "0.1227"
$ws.Range("E50").Value = "  +0.08%  "

# Row 51
$ws.Range("D51").Value = "0.2541"
$ws.Range("E51").Value = "  +3.74%  "

# Restore default styling (removes the temporary Text number format,
# leaving cells without an explicit style, as in the source file).
$dataRange.Style = "Normal"
